# Hardware.xlsx upload update
# - Renames "Lighting" -> "Electrical"
# - Adds a bold "Product" / "Stock" header row to every sheet
# - Appends new stock rows to Electrical and Plumbing
# - Adds two brand new sheets, "Flooring" and "Lumber", with their own data
# - Restores per-sheet selections / active tab to match the saved session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Grab/rename existing sheets
# ---------------------------------------------------------------------
$wsHardware   = $wb.Worksheets.Item("Hardware")
$wsElectrical = $wb.Worksheets.Item("Lighting")
$wsElectrical.Name = "Electrical"
$wsPlumbing   = $wb.Worksheets.Item("Plumbing")

# ---------------------------------------------------------------------
# Add the two new sheets at the end, in order: Flooring, then Lumber
# ---------------------------------------------------------------------
$wsFlooring = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsFlooring.Name = "Flooring"
$wsLumber = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsLumber.Name = "Lumber"

# ---------------------------------------------------------------------
# Hardware: add header row (data rows are unchanged)
# ---------------------------------------------------------------------
$wsHardware.Range("A1").Value = "Product"
$wsHardware.Range("B1").Value = "Stock"
$wsHardware.Range("A1:B1").Font.Bold = $true

# ---------------------------------------------------------------------
# Electrical: add header row + 3 new stock rows
# ---------------------------------------------------------------------
$wsElectrical.Range("A1").Value = "Product"
$wsElectrical.Range("B1").Value = "Stock"
$wsElectrical.Range("A1:B1").Font.Bold = $true

$wsElectrical.Range("A5").Value = "Breakers"
$wsElectrical.Range("B5").Value = 100
$wsElectrical.Range("A6").Value = "Extension Cords"
$wsElectrical.Range("B6").Value = 50
$wsElectrical.Range("A7").Value = "Fuses"
$wsElectrical.Range("B7").Value = 300

# ---------------------------------------------------------------------
# Plumbing: add header row + 2 new stock rows
# ---------------------------------------------------------------------
$wsPlumbing.Range("A1").Value = "Product"
$wsPlumbing.Range("B1").Value = "Stock"
$wsPlumbing.Range("A1:B1").Font.Bold = $true

$wsPlumbing.Range("A6").Value = "Caulking"
$wsPlumbing.Range("B6").Value = 200
$wsPlumbing.Range("A7").Value = "Valves"
$wsPlumbing.Range("B7").Value = 250

# ---------------------------------------------------------------------
# Flooring: header row + full data set
# ---------------------------------------------------------------------
$wsFlooring.Range("A1").Value = "Product"
$wsFlooring.Range("B1").Value = "Stock"
$wsFlooring.Range("A1:B1").Font.Bold = $true

$wsFlooring.Range("A2").Value = "Marble Tile"
$wsFlooring.Range("B2").Value = 500
$wsFlooring.Range("A3").Value = "Carpet (ft.)"
$wsFlooring.Range("B3").Value = 500
$wsFlooring.Range("A4").Value = "Vacuum"
$wsFlooring.Range("B4").Value = 25
$wsFlooring.Range("A5").Value = "Floor Finish"
$wsFlooring.Range("B5").Value = 60
$wsFlooring.Range("A6").Value = "Versabond"
$wsFlooring.Range("B6").Value = 100
$wsFlooring.Range("A7").Value = "Wood Tile"
$wsFlooring.Range("B7").Value = 500

$wsFlooring.Columns.Item(1).ColumnWidth = 9.666666666666666

# ---------------------------------------------------------------------
# Lumber: header row + full data set
# ---------------------------------------------------------------------
$wsLumber.Range("A1").Value = "Product"
$wsLumber.Range("B1").Value = "Stock"
$wsLumber.Range("A1:B1").Font.Bold = $true

$wsLumber.Range("A2").Value = "2x4x8"
$wsLumber.Range("B2").Value = 200
$wsLumber.Range("A3").Value = "Sheetrock"
$wsLumber.Range("B3").Value = 150
$wsLumber.Range("A4").Value = "Chainlink"
$wsLumber.Range("B4").Value = 25
$wsLumber.Range("A5").Value = "Pickets"
$wsLumber.Range("B5").Value = 600
$wsLumber.Range("A6").Value = "Mailbox"
$wsLumber.Range("B6").Value = 30
$wsLumber.Range("A7").Value = "Concrete Mix"
$wsLumber.Range("B7").Value = 250

$wsLumber.Columns.Item(1).ColumnWidth = 10.833333333333334

# ---------------------------------------------------------------------
# Restore per-sheet selections; select Flooring last so it ends up the
# active tab, matching the saved workbook view.
# ---------------------------------------------------------------------
$wsHardware.Range("B1").Select()
$wsElectrical.Range("A9").Select()
$wsPlumbing.Range("C8").Select()
$wsLumber.Range("C9").Select()
$wsFlooring.Range("C4").Select()
